$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the R19:R23 and R27:R31 "crosstab_percent" formulas: the multiplier
# (4/80) must be applied inside the SQRT, not outside it.
$ws.Range("R19").Formula = "=SQRT((4/80)*SUM(N19:Q19))"
$ws.Range("R20:R23").Formula = "=SQRT((4/80)*SUM(N20:Q20))"
$ws.Range("R27").Formula = "=SQRT((4/80)*SUM(N27:Q27))"
$ws.Range("R28:R31").Formula = "=SQRT((4/80)*SUM(N28:Q28))"

# Move the active selection to reflect where the author left off reviewing.
$ws.Range("R28").Select()
